# Homework #3: add "Problem 4." and "Problem 5." sections (with answers)
# right before the trailing blank paragraph at the end of the document body.

$d = $word.ActiveDocument

# The document ends with: [...last picture...] [blank Courier-New paragraph]
# [blank "Normal" paragraph]. We want to insert the new content between the
# blank Courier-New paragraph and the final blank paragraph, leaving both of
# those untouched.
$lastPara = $d.Paragraphs.Last
$anchorPara = $lastPara.Previous()

# Create a fresh, empty paragraph right after the anchor. This new paragraph
# becomes our insertion target; Range.InsertXML() *replaces* the contents of
# whatever range it is called on, so operating on this disposable paragraph's
# full range (not a collapsed point) lets the supplied OOXML become brand new
# sibling paragraphs without merging into -- or stealing the formatting of --
# any pre-existing paragraph.
$anchorPara.Range.InsertParagraphAfter() | Out-Null

$count = $d.Paragraphs.Count
$placeholder = $d.Paragraphs.Item($count - 1)
$target = $d.Range($placeholder.Range.Start, $placeholder.Range.End)

$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Problem 4.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">4. (20%) What is light field imaging? Describe the advantages and limitations of light field imaging. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Answer:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>A fundamental definition of a light field would be the parameterization of the flow of light thru an empty region of 3D space. Now</w:t></w:r><w:r><w:t>, there are many ways to represent these parameters. In the most redundant of systems you can use a 7D system, which would include (x, y, z, theta(angle), alpha(angle), t, and lambda(wavelength). Th</w:t></w:r><w:r><w:t xml:space="preserve">ere are other systems, in the </w:t></w:r><w:r><w:t>paper “</w:t></w:r><w:r><w:t>Light Field and Computational Imaging”, Levoy references Moon and Hanrahan which annotate a 4D system as a light field.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">There are </w:t></w:r><w:r><w:t>some</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>limitations</w:t></w:r><w:r><w:t xml:space="preserve"> of light-field imaging. </w:t></w:r><w:r><w:t>One</w:t></w:r><w:r><w:t xml:space="preserve"> disadvantage is that occlusion prevents information </w:t></w:r><w:r><w:t>of concave objects. Also, if time is part of your parameters there are natural limitations of illumination of dynamic scenes.</w:t></w:r><w:r><w:t xml:space="preserve"> Although not a strict disadvantage</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> as the</w:t></w:r><w:r><w:t>re are ways to compensate, but geometrical complex objects present challenges to light field processing as multiple reflections, refractions and scattering affect the light ray. This last area is an active area of research.</w:t></w:r><w:r><w:t xml:space="preserve"> Perhaps the biggest disadvantage</w:t></w:r><w:r><w:t xml:space="preserve">s are economics and </w:t></w:r><w:r><w:t xml:space="preserve">adoption of the technology. Because of the ascendancy of CMOS cameras in ubiquitous smart phones like the iPhone, the public does not see much extra value with plenoptic functions and light field systems. The widespread adoption of light fields would have to add salient extra benefits not found in today’s cameras. However, to take advantage of that would require a technology increase in bandwidth, sensors, </w:t></w:r><w:r><w:t>memory,</w:t></w:r><w:r><w:t xml:space="preserve"> and computer speed many times greater than today to make light fields camera</w:t></w:r><w:r><w:t xml:space="preserve"> applications practical.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">On the other hand, light field imaging provides new ways of capturing images and brings some </w:t></w:r><w:r><w:t xml:space="preserve">benefits that traditional photography and imaging do not have. </w:t></w:r><w:r><w:t>One advantage of light-field imaging is that with one shot you can provide “plenoptic” information</w:t></w:r><w:r><w:t xml:space="preserve"> of angles which would allow</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> with computational imaging</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>the reconstruction of focusing on different regions of the image.</w:t></w:r><w:r><w:t xml:space="preserve"> Closely related to this advantage is that you do not have to decide </w:t></w:r><w:r><w:t xml:space="preserve">at the moment </w:t></w:r><w:r><w:t>of image capture the focus as “post-processing” can extract different depths of field.</w:t></w:r><w:r><w:t xml:space="preserve"> Another advantage is that there are at least two different ways</w:t></w:r><w:r><w:t xml:space="preserve"> to capture light fields in an inexpensive way: sequential capture and spatial multiplexing. </w:t></w:r><w:r><w:t xml:space="preserve">Spatial multiplexing uses a </w:t></w:r><w:r><w:t>“robotic arm”</w:t></w:r><w:r><w:t xml:space="preserve"> to position lamps and cameras to capture the plenoptic information, while spatial multiplexing makes use of micro-lenses.</w:t></w:r><w:r><w:t xml:space="preserve"> Finally, because light field provides more information than just 2-D, views and angles that were not even captured or sampled with the given light field system, a light field rendering can generate new images, a form of interpolation that provides new custom views and presentations of the original image.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/></w:pPr><w:r><w:t>Problem 5.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">5. (20%) Draw a diagram of a typical light-field camera, indicate the key </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>components,</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and explain their functions. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Default"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Answer:</w:t></w:r></w:p>'

$target.InsertXML($newContentXml)
